$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 518, which shifts existing rows 518-556 down to 519-557
$ws.Rows("518:518").Insert()

# Populate the newly inserted row 518 with the new weekly data point.
# Most columns stay the same as the surrounding rows for this subset (market/category metadata).
$ws.Cells.Item(518, 1).Value = 4
$ws.Cells.Item(518, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(518, 3).Value = "Los Lagos"
$ws.Cells.Item(518, 4).Value = 45265
$ws.Cells.Item(518, 5).Value = 10
$ws.Cells.Item(518, 6).Value = 100112037
$ws.Cells.Item(518, 7).Value = "Cebollín"
$ws.Cells.Item(518, 8).Value = "Sin especificar"
$ws.Cells.Item(518, 9).Value = "Primera"
$ws.Cells.Item(518, 10).Value = 180
$ws.Cells.Item(518, 11).Value = 6000
$ws.Cells.Item(518, 12).Value = 6000
$ws.Cells.Item(518, 13).Value = 6000
$ws.Cells.Item(518, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(518, 15).Value = "Región Metropolitana"
$ws.Cells.Item(518, 16).Value = 167
$ws.Cells.Item(518, 17).Value = 36
$ws.Cells.Item(518, 18).Value = "Hortaliza"

# Fix the Origen value swap caused by the row shift: what is now row 541 (previously row 540)
# should read "Provincia de Chacabuco", and what is now row 540 (previously row 539) should read
# "Región Metropolitana".
$ws.Cells.Item(540, 15).Value = "Región Metropolitana"
$ws.Cells.Item(541, 15).Value = "Provincia de Chacabuco"
